# Fill in the half-term summary text for Navi, Geometry, and Arithmetic
# subjects on the certificate (table content cells were empty placeholders).
$d = $word.ActiveDocument

# Table 2 = "נביא" (Navi/Prophets) row
$d.Tables.Item(2).Cell(1, 2).Range.Text = "במחצית זאת למדנו את ספר שמואל, עם מפרשים והרחבנו בהרבה מדרשים,
הייתה אוירה קדושה,והנאה בלימוד.
אהובה את מיוחדת, בהצלחה!"

# Table 3 = "הנדסה" (Geometry) row
$d.Tables.Item(3).Cell(1, 2).Range.Text = "למדנו במחצית זאת על סוגי המצולעים השונים. כגון ריבוע, משולש, משושה, מעוין וכו'. למדנו כיצד מזהים כל מצולע ומה תכונותיו.
אהובה, את תלמידה מדהימה! המשיכי להצליח בדרכך!"

# Table 4 = "חשבון" (Arithmetic) row
$d.Tables.Item(4).Cell(1, 2).Range.Text = "במחצית זאת חזרנו את פעולות חשבון, חיבור, חיסור, כפל וחילוק,
התמקדנו בעיקר על לוח הכפל, פיתחנו שיטות לימוד רציניות,ולמדנו איך לזכור דברים בע`"פ.
אהובה את ילדה מקסימה, יש לך ראש חזק, שיהיה בהצלחה!"
